$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The branch's client/employee counts were reset to 0 (C2: Num_Clients, D2: Num_Employees)
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

# Move/restore the cell selection to D3 (was E7)
$ws.Range("D3").Select()

# Reflect the updated Excel window size/position from the saved view state
$win = $excel.ActiveWindow
$win.Left = 7245
$win.Top = 525
$win.Width = 21600
$win.Height = 11385
